$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.021.84"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "1.824.21"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").Value = "311.52"

$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").Value = "0.4694"
$ws.Range("E7").Value = "  +1.41%  "

$ws.Range("D8").Value = "0.3665"
$ws.Range("E8").Value = "  -0.78%  "

$ws.Range("D9").Value = "0.07358"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "0.8752"
$ws.Range("E10").Value = "  +0.72%  "

$ws.Range("D11").Value = "20.31"
$ws.Range("E11").Value = "  -0.34%  "

$ws.Range("D12").Value = "1.819.92"
$ws.Range("E12").Value = "  -4.53%  "

$ws.Range("D13").Value = "0.07315"
$ws.Range("E13").Value = "  +3.42%  "

$ws.Range("D14").Value = "5.431"
$ws.Range("E14").Value = "  +1.69%  "

$ws.Range("D15").Value = "6.523"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").Value = "91.88"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "0.000008743"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("D21").Value = "27.036.37"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "5.295"
$ws.Range("E22").Value = "  -0.62%  "

$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("D24").Value = "2.063.36"
$ws.Range("E24").Value = "  -2.14%  "

$ws.Range("D25").Value = "1.896"
$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "151.45"
$ws.Range("E26").Value = "  -0.29%  "

$ws.Range("E27").Value = "  +0.41%  "

$ws.Range("D28").Value = "2.147"
$ws.Range("E28").Value = "  +1.39%  "

$ws.Range("D29").Value = "5.237"
$ws.Range("E29").Value = "  -0.96%  "

$ws.Range("D30").Value = "116.71"
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").Value = "0.08894"
$ws.Range("E31").Value = "  +0.07%  "

$ws.Range("D32").Value = "0.7541"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "1.162"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").Value = "4.513"
$ws.Range("E34").Value = "  +1.27%  "

$ws.Range("D35").Value = "2.932"
$ws.Range("E35").Value = "  +0.47%  "

$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  +0.11%  "

$ws.Range("D38").Value = "0.05310"
$ws.Range("E38").Value = "  +1.10%  "

$ws.Range("D39").Value = "0.01952"
$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").Value = "2.978"
$ws.Range("E40").Value = "  +2.24%  "

$ws.Range("D41").Value = "7.221"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("D42").Value = "2.378"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").Value = "0.5309"
$ws.Range("E43").Value = "  -0.54%  "

$ws.Range("D44").Value = "0.1655"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").Value = "8.486"
$ws.Range("E45").Value = "  +0.83%  "

$ws.Range("D46").Value = "0.4902"
$ws.Range("E46").Value = "  -0.60%  "

$ws.Range("D47").Value = "10.42"
$ws.Range("E47").Value = "  +0.72%  "

$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("D49").Value = "1.665"
$ws.Range("E49").Value = "  -0.32%  "

$ws.Range("D50").Value = "103.18"
$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").Value = "0.06303"
$ws.Range("E51").Value = "  +0.61%  "
